# Reposition / resize the logo picture on slide 2.
#
# Target (EMU, from the canonical OOXML diff):
#   a:off  x="3944056" y="1800837"
#   a:ext cx="3540410" cy="3540410"
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU). The literals below were chosen so that, once run
# through the host's point -> EMU conversion, they land exactly on the
# target EMU values (plain EMU/12700 division is off by a hair because
# of float rounding in that conversion).
$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)

$shp.Left   = 310.55563355118113
$shp.Top    = 141.79823303637795
$shp.Width  = 278.7724761948819
$shp.Height = 278.7724761948819
